$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.257390022277832
$ws.Range("B1").Value = 2.313507795333862
$ws.Range("C1").Value = 3.697669506072998
$ws.Range("D1").Value = 2.758620977401733
$ws.Range("E1").Value = 1.361572742462158
